$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.021.24'
$ws.Range("E2").Value = '  +7.88%  '
$ws.Range("D3").Value = '3.319.62'
$ws.Range("E3").Value = '  +2.84%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.04'
$ws.Range("E5").Value = '  +4.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '627.59'
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.412'
$ws.Range("E7").Value = '  +14.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.707'
$ws.Range("E8").Value = '  +5.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '3.304.54'
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.591'
$ws.Range("E11").Value = '  +5.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000266'
$ws.Range("E12").Value = '  +6.99%  '
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.88'
$ws.Range("E14").Value = '  +4.31%  '
$ws.Range("D15").Value = '3.935.76'
$ws.Range("E15").Value = '  +2.82%  '
$ws.Range("D16").Value = '92.344.95'
$ws.Range("E16").Value = '  +7.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.38'
$ws.Range("E17").Value = '  +2.90%  '
$ws.Range("D18").Value = '3.292.92'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.36'
$ws.Range("E19").Value = '  +12.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.25'
$ws.Range("E20").Value = '  +3.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '439.74'
$ws.Range("E21").Value = '  +3.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.96'
$ws.Range("E22").Value = '  +3.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.37'
$ws.Range("E23").Value = '  +2.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000191'
$ws.Range("E24").Value = '  +52.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.40'
$ws.Range("E25").Value = '  +8.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.49'
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("D27").Value = '3.521.55'
$ws.Range("E27").Value = '  +3.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.71'
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.180'
$ws.Range("E30").Value = '  +6.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.85'
$ws.Range("E32").Value = '  +2.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '569.04'
$ws.Range("E33").Value = '  +6.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.34'
$ws.Range("E34").Value = '  +7.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.71'
$ws.Range("E35").Value = '  +30.29%  '
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.93'
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.33'
$ws.Range("E37").Value = '  -3.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.87'
$ws.Range("E38").Value = '  +3.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.132'
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.46'
$ws.Range("E40").Value = '  +4.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.398'
$ws.Range("E42").Value = '  +3.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.01'
$ws.Range("E43").Value = '  +3.12%  '
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '151.01'
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '182.90'
$ws.Range("E46").Value = '  +3.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.130'
$ws.Range("E47").Value = '  +8.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.74'
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.29'
$ws.Range("E49").Value = '  +2.85%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.636'
$ws.Range("E50").Value = '  +4.13%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.19'
$ws.Range("E51").Value = '  +1.33%  '
